$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.219.03'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.661.68'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.40'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5215'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2644'
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06285'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.83'
$ws.Range('E10').Value = '  -3.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07778'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.478'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.645.96'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = '1.889.21'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5462'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = '0.0₅8162'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.97'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '26.219.91'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.604'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.98'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.03'
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.012'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.06'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1230'
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.275'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.19'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.425'
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05966'
$ws.Range('E30').Value = '  -3.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.276'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.544'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.274'
$ws.Range('E33').Value = '  -3.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.581'
$ws.Range('E34').Value = '  -5.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9603'
$ws.Range('E35').Value = '  -4.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.416'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5686'
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01599'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.976'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8502'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.56'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.003.14'
$ws.Range('E44').Value = '  -8.01%  '
$ws.Range('D45').Value = '1.803.97'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.61'
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.026'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05150'
